$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Cat3" column (J) header
$ws.Range("J1").Value = "Cat3"

# Fill in the Cat3 category values for rows 2-32
$cat3Values = @("A","S","F","F","A","S","Z","X","C","Z","S","F","C","V","Z","X","C","V","Z","S","A","C","F","A","S","A","A","A","F","Z","Z")

for ($i = 0; $i -lt $cat3Values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $cat3Values[$i]
}

# Update the active selection to reflect the new column being worked on
$null = $ws.Range("J4").Select()
